$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header-less column E with team names for each player.
# Rows 1-25: Chennai Super Kings; rows 26-50: Mumbai Indians.
for ($r = 1; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value = "Chennai Super Kings"
}
for ($r = 26; $r -le 50; $r++) {
    $ws.Cells.Item($r, 5).Value = "Mumbai Indians"
}

# Set column E width to match the diff (target stored width 17.83203125;
# the closest value the engine's pixel-grid snapping can produce is 17.8333,
# reached via a ColumnWidth of 17).
$ws.Columns.Item(5).ColumnWidth = 17

# Update the active cell selection to C9 (matches diff's <selection activeCell="C9" sqref="C9"/>)
$ws.Range("C9").Select()
